# Auto-generated edit script applying market-data refresh to computed columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 663.26666
$ws.Range("I28").Value = 692.2308
$ws.Range("J28").Value = 475
$ws.Range("K28").Value = 692.2308
$ws.Range("L28").Value = 475
$ws.Range("M28").Value = -207.2308
$ws.Range("N28").Value = -1445
$ws.Range("H33").Value = 338.05884
$ws.Range("I33").Value = 296.8125
$ws.Range("K33").Value = 296.8125
$ws.Range("M33").Value = -67.8125
$ws.Range("H51").Value = 316599.66
$ws.Range("J51").Value = 274899.5
$ws.Range("L51").Value = 274899.5
$ws.Range("N51").Value = -275867.5
$ws.Range("H58").Value = 2591.5715
$ws.Range("I58").Value = 785.75
$ws.Range("J58").Value = 4999.3335
$ws.Range("K58").Value = 2357.25
$ws.Range("L58").Value = 14998.0005
$ws.Range("M58").Value = -2207.25
$ws.Range("N58").Value = -15298.0005
$ws.Range("H105").Value = 45333.332
$ws.Range("J105").Value = 45333.332
$ws.Range("L105").Value = 45333.332
$ws.Range("N105").Value = -52321.332
$ws.Range("H107").Value = 501.625
$ws.Range("I107").Value = 501.625
$ws.Range("K107").Value = 501.625
$ws.Range("M107").Value = 1418.375
$ws.Range("H127").Value = 500
$ws.Range("I127").Value = 500
$ws.Range("K127").Value = 1500
$ws.Range("M127").Value = 3460
$ws.Range("H138").Value = 2666.5789
$ws.Range("I138").Value = 1299
$ws.Range("K138").Value = 3897
$ws.Range("M138").Value = 1243

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1234.0555
$ws.Range("I97").Value = 758.0714
$ws.Range("K97").Value = 758.0714
$ws.Range("M97").Value = -262.0714
$ws.Range("H110").Value = 788.4545000000001
$ws.Range("I110").Value = 936.75
$ws.Range("K110").Value = 936.75
$ws.Range("M110").Value = 1108.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1051.5
$ws.Range("I5").Value = 1051.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1051.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -938.5
$ws.Range("N5").ClearContents()
$ws.Range("H20").Value = 1716.8334
$ws.Range("I20").Value = 1716.8334
$ws.Range("K20").Value = 1716.8334
$ws.Range("M20").Value = -1469.8334
$ws.Range("H26").Value = 25289.25
$ws.Range("I26").Value = 25289.25
$ws.Range("K26").Value = 25289.25
$ws.Range("M26").Value = -24997.25
$ws.Range("H107").Value = 4057.6191
$ws.Range("I107").Value = 4057.6191
$ws.Range("K107").Value = 4057.6191
$ws.Range("M107").Value = -2137.6191

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H31").Value = 1900
$ws.Range("I31").Value = 1900
$ws.Range("K31").Value = 1900
$ws.Range("M31").Value = -1605
$ws.Range("H34").Value = 1900
$ws.Range("I34").Value = 1900
$ws.Range("K34").Value = 1900
$ws.Range("M34").Value = -1698
$ws.Range("H99").Value = 2578.6
$ws.Range("I99").Value = 2578.6
$ws.Range("K99").Value = 2578.6
$ws.Range("M99").Value = -1080.6
$ws.Range("H107").Value = 640.125
$ws.Range("I107").Value = 603.5
$ws.Range("K107").Value = 603.5
$ws.Range("M107").Value = 1316.5
$ws.Range("H126").Value = 2578.6
$ws.Range("I126").Value = 2578.6
$ws.Range("K126").Value = 7735.799999999999
$ws.Range("M126").Value = -5265.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1133
$ws.Range("J5").Value = 2000
$ws.Range("L5").Value = 6000
$ws.Range("N5").Value = -6224
$ws.Range("H7").Value = 370.6
$ws.Range("I7").Value = 337.75
$ws.Range("J7").Value = 502
$ws.Range("K7").Value = 1013.25
$ws.Range("L7").Value = 1506
$ws.Range("M7").Value = -901.25
$ws.Range("N7").Value = -1730
$ws.Range("H120").Value = 25000
$ws.Range("I120").Value = 25000
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 75000
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -70162
$ws.Range("N120").ClearContents()
$ws.Range("H129").Value = 230
$ws.Range("I129").Value = 230
$ws.Range("K129").Value = 690
$ws.Range("M129").Value = 4310
$ws.Range("H135").Value = 1133
$ws.Range("J135").Value = 2000
$ws.Range("L135").Value = 18000
$ws.Range("N135").Value = -23070
$ws.Range("H138").Value = 4021.6
$ws.Range("I138").Value = 4021.6
$ws.Range("K138").Value = 12064.8
$ws.Range("M138").Value = -6924.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 4
$ws.Range("K5").Value = 4
$ws.Range("M5").Value = 108
$ws.Range("H48").Value = 8000
$ws.Range("I48").Value = 8000
$ws.Range("K48").Value = 8000
$ws.Range("M48").Value = -7515
$ws.Range("H113").Value = 1497.5
$ws.Range("I113").Value = 1497.5
$ws.Range("K113").Value = 1497.5
$ws.Range("M113").Value = 672.5
$ws.Range("H124").Value = 19000
$ws.Range("I124").Value = 19000
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 19000
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = -14090
$ws.Range("N124").ClearContents()
$ws.Range("H133").Value = 145000
$ws.Range("J133").Value = 145000
$ws.Range("L133").Value = 145000
$ws.Range("N133").Value = -155120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 664.8
$ws.Range("I55").Value = 474.5
$ws.Range("J55").Value = 791.6667
$ws.Range("K55").Value = 474.5
$ws.Range("L55").Value = 791.6667
$ws.Range("M55").Value = -301.5
$ws.Range("N55").Value = -1137.6667
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 500
$ws.Range("J4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -726
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H54").Value = 28496.25
$ws.Range("J54").Value = 28496.25
$ws.Range("L54").Value = 28496.25
$ws.Range("N54").Value = -29536.25
$ws.Range("H132").Value = 3212
$ws.Range("I132").Value = 3212
$ws.Range("K132").Value = 9636
$ws.Range("M132").Value = -7106
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
